{"js": "// Office.js (Word JavaScript API) edit script.\n// Body of: async (context) => { ... }\n//\n// Changes applied (per the diff):\n// 1. Insert a new Heading-1 (\"Titre1\") paragraph reading \"This is some text\"\n//    right after the \"Intro of part 2\" paragraph (and before \"Subpart 1\").\n// 2. Merge the two runs \"Content of subpart \" + \"2\" into a single run\n//    \"Content of subpart 2\" inside the paragraph that follows \"Little title 2\".\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nlet introParagraph = null;\nlet subpart2Paragraph = null;\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const paragraph = paragraphs.items[i];\n  const text = paragraph.text;\n  if (text === \"Intro of part 2\") {\n    introParagraph = paragraph;\n  } else if (text === \"Content of subpart 2\") {\n    subpart2Paragraph = paragraph;\n  }\n}\n\nif (!introParagraph) {\n  throw new Error('Could not find paragraph \"Intro of part 2\"');\n}\n\n// 1. Insert the new heading paragraph after \"Intro of part 2\".\nconst newParagraph = introParagraph.insertParagraph(\"This is some text\", Word.InsertLocation.after);\nnewParagraph.styleBuiltIn = Word.Style.heading1;\n\n// 2. Replace the split-run paragraph's text with a single merged run.\nif (subpart2Paragraph) {\n  subpart2Paragraph.insertText(\"Content of subpart 2\", Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n# $word.ActiveDocument is already open as $d below.\n#\n# Changes applied (per the diff):\n# 1. Merge the two runs \"Content of subpart \" + \"2\" into a single run\n#    \"Content of subpart 2\" (paragraph that follows \"Little title 2\").\n#    Done FIRST so paragraph indices/offsets used for step 2 are not\n#    disturbed by the later structural insertion.\n# 2. Insert a new Heading-1 (\"Titre1\") paragraph reading \"This is some text\"\n#    right after the \"Intro of part 2\" paragraph (and before \"Subpart 1\").\n\n$d = $word.ActiveDocument\n\n# --- Step 1: merge the split run \"Content of subpart \" + \"2\" ---\n$subpart2Para = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text -eq (\"Content of subpart 2\" + [char]13)) {\n        $subpart2Para = $p\n    }\n}\n\nif ($subpart2Para -ne $null) {\n    $start = $subpart2Para.Range.Start\n    # \"Content of subpart \" is 19 characters, \"2\" is 1 character.\n    $firstRunRange = $d.Range($start, $start + 19)\n    $secondRunRange = $d.Range($start + 19, $start + 20)\n    $secondRunRange.Text = \"\"\n    $firstRunRange.Text = \"Content of subpart 2\"\n}\n\n# --- Step 2: insert the new heading paragraph after \"Intro of part 2\" ---\n# Re-query for the paragraph fresh (do not reuse a reference captured\n# before the step-1 mutation) so its Range reflects current offsets.\n$introPara = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text -eq (\"Intro of part 2\" + [char]13)) {\n        $introPara = $p\n    }\n}\n\nif ($introPara -ne $null) {\n    $rng = $introPara.Range\n    $rng.Collapse(0)\n    $rng.InsertParagraphAfter()\n    $newPara = $introPara.Next()\n    $newPara.Range.Text = \"This is some text\"\n    $newPara.Style = \"Titre1\"\n}\n"}
